# Apply the "years 1-3 for aggressive" edit to the management strategies sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the shared text describing the Aggressive Strategy's mechanical + chemical
#    treatment schedule: "years 1-2" -> "years 1-3". All cells H5:H16 share this string,
#    so updating the text (in place) propagates to every row automatically.
$oldText = "Mechanical + chemical treatment in years 1-2 every 5 years; cover crop re-seeded after each treatment"
$newText = "Mechanical + chemical treatment in years 1-3 every 5 years; cover crop re-seeded after each treatment"

for ($r = 5; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 8)  # column H
    if ($cell.Value() -eq $oldText) {
        $cell.Value = $newText
    }
}

# 2. Increase the row heights for the data rows (5-15 from 48 to 68, row 16 from 49 to 69)
#    to accommodate the longer wrapped text.
for ($r = 5; $r -le 15; $r++) {
    $ws.Rows.Item($r).RowHeight = 68
}
$ws.Rows.Item(16).RowHeight = 69

# 3. Scroll the sheet view so row 7 is the top-left visible cell.
$ws.Application.ActiveWindow.ScrollRow = 7
